$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

# Add a new "metadata" sheet right after the existing "data" sheet.
$ws = $wb.Worksheets.Add($null, $data)
$ws.Name = "metadata"

# Match the header/"key" cell formatting used on the data sheet (bold,
# centered, thin border) by copying it over from an already-styled cell.
$data.Range("F1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Data row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Peroxisomal Disorders"
$ws.Range("C2").Value = 155
$ws.Range("D2").Value = "'0.23"
$ws.Range("E2").Value = "2021-07-27T23:11:55.003702Z"
$ws.Range("F2").Value = "2021-10-05 14:35:12.687380"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/155/?format=json"

# Refresh the per-row "time_taken" timestamps on the data sheet.
$times = @(
  "2021-10-05 14:35:12.691088",
  "2021-10-05 14:35:12.691096",
  "2021-10-05 14:35:12.691099",
  "2021-10-05 14:35:12.691102",
  "2021-10-05 14:35:12.691105",
  "2021-10-05 14:35:12.691107",
  "2021-10-05 14:35:12.691110",
  "2021-10-05 14:35:12.691112",
  "2021-10-05 14:35:12.691115",
  "2021-10-05 14:35:12.691118",
  "2021-10-05 14:35:12.691120",
  "2021-10-05 14:35:12.691123",
  "2021-10-05 14:35:12.691125",
  "2021-10-05 14:35:12.691128",
  "2021-10-05 14:35:12.691130",
  "2021-10-05 14:35:12.691133",
  "2021-10-05 14:35:12.691135",
  "2021-10-05 14:35:12.691138",
  "2021-10-05 14:35:12.691141",
  "2021-10-05 14:35:12.691143",
  "2021-10-05 14:35:12.691146",
  "2021-10-05 14:35:12.691148",
  "2021-10-05 14:35:12.691151",
  "2021-10-05 14:35:12.691153",
  "2021-10-05 14:35:12.691156",
  "2021-10-05 14:35:12.691158",
  "2021-10-05 14:35:12.691161",
  "2021-10-05 14:35:12.691163",
  "2021-10-05 14:35:12.691166",
  "2021-10-05 14:35:12.691168",
  "2021-10-05 14:35:12.691171",
  "2021-10-05 14:35:12.691173"
)

for ($i = 0; $i -lt $times.Length; $i++) {
  $row = $i + 2
  $data.Cells.Item($row, 6).Value = $times[$i]
}
